# Apply the VolSkill.xlsx update:
#  - Drop the stray "VolTeer.Vol.tblVolSkill" access-string cell (old B1)
#  - Promote the header row (Table/SkillID/VolID/Query) up to row 1
#  - Shift the three data rows up, replace the two GUID values in the
#    last data row (VolID column) and regenerate the INSERT-INTO helper
#    formula on every data row
#  - Drop the now-superfluous 5th row entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing row (was row 5) so only 3 data rows remain.
$ws.Rows.Item(5).Delete()

# Row 1: table-name / query header, moved up from row 2.
$ws.Range("A1").Value = "Table"
$ws.Range("B1").Value = "SkillID"
$ws.Range("C1").Value = "VolID"
$ws.Range("D1").Value = "Query"

# Row 2 data + formula.
$ws.Range("A2").Value = "Vol.tblVolSkill"
$ws.Range("B2").Value = "153B9006-1488-4CA1-950A-32E28F7B699D"
$ws.Range("C2").Value = "dba53101-f9b2-4dc0-85e7-11d472fd99cd"
$ws.Range("D2").Formula = '=((((((((("INSERT INTO " & A2) &" (") & $B$1) & ",") & $C$1) & ") VALUES(''") & RC[-2]) &  "'',''") & RC[-1]) & "'')"'

# Row 3 data + formula.
$ws.Range("A3").Value = "Vol.tblVolSkill"
$ws.Range("B3").Value = "C87F23E9-8F8C-406D-9FBF-E15043179F09"
$ws.Range("C3").Value = "589178b4-aa4c-4276-a516-9460fa7714d3"
$ws.Range("D3").Formula = '=((((((((("INSERT INTO " & A3) &" (") & $B$1) & ",") & $C$1) & ") VALUES(''") & RC[-2]) &  "'',''") & RC[-1]) & "'')"'

# Row 4 data + formula.
$ws.Range("A4").Value = "Vol.tblVolSkill"
$ws.Range("B4").Value = "990B091D-6A0D-4543-B766-C37B684F8081"
$ws.Range("C4").Value = "293fe520-7e35-444a-8955-f02a911fed1c"
$ws.Range("D4").Formula = '=((((((((("INSERT INTO " & A4) &" (") & $B$1) & ",") & $C$1) & ") VALUES(''") & RC[-2]) &  "'',''") & RC[-1]) & "'')"'
